# Applies the experimental_data.xlsx edits described by the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet "Pre" ---------------------------------------------------------
$pre = $wb.Worksheets.Item("Pre")
$pre.Activate()
$pre.Range("B3:C4").Select()

# --- Sheet "Main" ---------------------------------------------------------
$main = $wb.Worksheets.Item("Main")
$main.Activate()

# Row 6
$main.Cells.Item(6, 2).Value = 310
$main.Cells.Item(6, 3).Value = 0.222

# Row 7
$main.Cells.Item(7, 2).Value = 315
$main.Cells.Item(7, 3).Value = 0.22500000000000001

# Row 9
$main.Cells.Item(9, 2).Value = 269
$main.Cells.Item(9, 3).Value = 0.2152

# Row 10
$main.Cells.Item(10, 2).Value = 263
$main.Cells.Item(10, 3).Value = 0.2162

# Row 11
$main.Cells.Item(11, 2).Value = 187
$main.Cells.Item(11, 3).Value = 0.20649999999999999

# Row 12
$main.Cells.Item(12, 2).Value = 269
$main.Cells.Item(12, 3).Value = 0.2152

# Row 13
$main.Cells.Item(13, 2).Value = 263
$main.Cells.Item(13, 3).Value = 0.2162

$main.Range("D20").Select()
